$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Row 6
$ws.Range("B6").Value = 0.7983299016243424
$ws.Range("C6").Value = 0.809081909970925
$ws.Range("D6").Value = 0.7983299016243424
$ws.Range("E6").Value = 0.7992955421404832
$ws.Range("F6").Value = 0.8154884465797301
$ws.Range("G6").Value = 0.8260992926735387
$ws.Range("H6").Value = 0.8154884465797301
$ws.Range("I6").Value = 0.8166684721988402
$ws.Range("J6").Value = 0.7274307938686799
$ws.Range("K6").Value = 0.7374723225254568
$ws.Range("L6").Value = 0.7274307938686799
$ws.Range("M6").Value = 0.7266662225643733
$ws.Range("N6").Value = 0.8111873713109128
$ws.Range("O6").Value = 0.8224609245887479
$ws.Range("P6").Value = 0.8111873713109128
$ws.Range("Q6").Value = 0.8123268668058288
$ws.Range("R6").Value = 0.8069778082818578
$ws.Range("S6").Value = 0.8170258830515044
$ws.Range("T6").Value = 0.8069778082818578
$ws.Range("U6").Value = 0.8086470213208156
$ws.Range("V6").Value = 0.8177076183939602
$ws.Range("W6").Value = 0.8209627026737867
$ws.Range("X6").Value = 0.8177076183939602
$ws.Range("Y6").Value = 0.8165292199912176

# Row 7
$ws.Range("B7").Value = 0.8411805078929306
$ws.Range("C7").Value = 0.8465748250243991
$ws.Range("D7").Value = 0.8411805078929306
$ws.Range("E7").Value = 0.8419413357459821
$ws.Range("F7").Value = 0.8605811027224892
$ws.Range("G7").Value = 0.8628371299062879
$ws.Range("H7").Value = 0.8605811027224892
$ws.Range("I7").Value = 0.8596024325052311
$ws.Range("J7").Value = 0.8304964539007094
$ws.Range("K7").Value = 0.8431882219375659
$ws.Range("L7").Value = 0.8304964539007094
$ws.Range("M7").Value = 0.8316139843172679
$ws.Range("N7").Value = 0.8433768016472204
$ws.Range("O7").Value = 0.8485605807079966
$ws.Range("P7").Value = 0.8433768016472204
$ws.Range("Q7").Value = 0.8435653735638207
$ws.Range("R7").Value = 0.8755204758636468
$ws.Range("S7").Value = 0.8798494501818164
$ws.Range("T7").Value = 0.8755204758636468
$ws.Range("U7").Value = 0.8750760578550937
$ws.Range("V7").Value = 0.8455273392816289
$ws.Range("W7").Value = 0.8503573639297519
$ws.Range("X7").Value = 0.8455273392816289
$ws.Range("Y7").Value = 0.8457190001481187
